$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.957.02'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.61%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.820.87'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.10%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.49'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '113.81'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.42%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.560'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.13%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  +6.74%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.53'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.82%  '

$ws.Range('E11').Value = '  -0.66%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0846'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.58%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.97'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.38%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.86'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.46%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.259.00'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.96%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.972'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.53%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.822.78'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.46%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.956.42'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.68%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.37'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +10.08%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.86%  '

$ws.Range('E21').Value = '  +4.34%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0979'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.84%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.67'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.47%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.54'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.61%  '

$ws.Range('E25').Value = '  +1.87%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.31'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.39%  '

$ws.Range('E27').Value = '  -0.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.164'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.30%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.56'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.94%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.65'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +11.77%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.29'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.70%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '52.90'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.16%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.20'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.82%  '

$ws.Range('E34').Value = '  +9.48%  '

$ws.Range('E35').Value = '  +3.58%  '

$ws.Range('E36').Value = '  +2.23%  '

$ws.Range('E37').Value = '  -0.17%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.98'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.44%  '

$ws.Range('E39').Value = '  +2.91%  '

$ws.Range('E40').Value = '  +3.62%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.58'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.57%  '

$ws.Range('E42').Value = '  +2.14%  '

$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.24'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.19%  '

$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '121.20'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.76%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.29'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.26%  '

$ws.Range('E46').Value = '  +8.83%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.141.12'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.76%  '

$ws.Range('E48').Value = '  +8.90%  '

$ws.Range('E49').Value = '  +12.43%  '

$ws.Range('B50').Value = 'TheGraph'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.224'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +18.80%  '

$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0322'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +16.07%  '
